# Add new flower entries to the "Plant Level Tracking" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plant Level Tracking")

# --- Append new plant rows (33-41), alphabetically, with Plant / Artist / Stages columns ---
$ws.Range("A33").Value = "aloe_vera"
$ws.Range("B33").Value = "Diconcilio"
$ws.Range("F33").Value = 7
$ws.Range("F33").NumberFormat = "0"

$ws.Range("A34").Value = "cherry_blossom"
$ws.Range("B34").Value = "Diconcilio"
$ws.Range("F34").Value = 7
$ws.Range("F34").NumberFormat = "0"

$ws.Range("A35").Value = "hyacinth"
$ws.Range("B35").Value = "Diconcilio"
$ws.Range("F35").Value = 7
$ws.Range("F35").NumberFormat = "0"

$ws.Range("A36").Value = "hydrangea"
$ws.Range("B36").Value = "Diconcilio"
$ws.Range("F36").Value = 7
$ws.Range("F36").NumberFormat = "0"

$ws.Range("A37").Value = "lily_of_the_valley"
$ws.Range("B37").Value = "Diconcilio"
$ws.Range("F37").Value = 7
$ws.Range("F37").NumberFormat = "0"

$ws.Range("A38").Value = "moon_flower"
$ws.Range("B38").Value = "Diconcilio"
$ws.Range("F38").Value = 7
$ws.Range("F38").NumberFormat = "0"

$ws.Range("A39").Value = "silent_princess"
$ws.Range("B39").Value = "Diconcilio"
$ws.Range("F39").Value = 7
$ws.Range("F39").NumberFormat = "0"

$ws.Range("A40").Value = "tulip"
$ws.Range("B40").Value = "Diconcilio"
$ws.Range("F40").Value = 7
$ws.Range("F40").NumberFormat = "0"

$ws.Range("A41").Value = "widow_tears"
$ws.Range("B41").Value = "Gerti"
$ws.Range("F41").Value = 6
$ws.Range("F41").NumberFormat = "0"

# --- Fill in the new entries for the level-6 / level-7 unlock columns (row 6 & 7) ---
$ws.Range("H6").Value = "aloe_vera"
$ws.Range("I6").Value = "cherry_blossom"
$ws.Range("L6").Value = "hyacinth"
$ws.Range("M6").Value = "hydrangea"

$ws.Range("I7").Value = "moon_flower"
$ws.Range("J7").Value = "lily_of_the_valley"
$ws.Range("K7").Value = "silent_princess"
$ws.Range("L7").Value = "tulip"
$ws.Range("M7").Value = "widow_tears"

# --- Update the active selection to match the author's final cursor position ---
$null = $ws.Range("J10").Select()
